$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 601 (existing rows 601-630 shift down to 603-632)
$ws.Range("A601:T602").EntireRow.Insert()

# --- New row 601: Crimpson Seedless ---
$ws.Cells.Item(601, 1).Value = 9
$ws.Cells.Item(601, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(601, 3).Value = "Metropolitana"
$ws.Cells.Item(601, 4).Value = 44706
$ws.Cells.Item(601, 5).Value = 13
$ws.Cells.Item(601, 6).Value = "Fruta"
$ws.Cells.Item(601, 7).Value = 100109
$ws.Cells.Item(601, 8).Value = "Uva"
$ws.Cells.Item(601, 9).Value = 100109001
$ws.Cells.Item(601, 10).Value = "Uva"
$ws.Cells.Item(601, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(601, 12).Value = "Primera"
$ws.Cells.Item(601, 13).Value = 330
$ws.Cells.Item(601, 14).Value = 8000
$ws.Cells.Item(601, 15).Value = 8000
$ws.Cells.Item(601, 16).Value = 8000
$ws.Cells.Item(601, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(601, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(601, 19).Value = 444
$ws.Cells.Item(601, 20).Value = 18

# --- New row 602: Red Globe ---
$ws.Cells.Item(602, 1).Value = 9
$ws.Cells.Item(602, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(602, 3).Value = "Metropolitana"
$ws.Cells.Item(602, 4).Value = 44706
$ws.Cells.Item(602, 5).Value = 13
$ws.Cells.Item(602, 6).Value = "Fruta"
$ws.Cells.Item(602, 7).Value = 100109
$ws.Cells.Item(602, 8).Value = "Uva"
$ws.Cells.Item(602, 9).Value = 100109001
$ws.Cells.Item(602, 10).Value = "Uva"
$ws.Cells.Item(602, 11).Value = "Red Globe"
$ws.Cells.Item(602, 12).Value = "Primera"
$ws.Cells.Item(602, 13).Value = 350
$ws.Cells.Item(602, 14).Value = 9000
$ws.Cells.Item(602, 15).Value = 9000
$ws.Cells.Item(602, 16).Value = 9000
$ws.Cells.Item(602, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(602, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(602, 19).Value = 500
$ws.Cells.Item(602, 20).Value = 18
